$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update of "Hortaliza, Femacal de La Calera - Espárragos" dataset:
# per-row Fecha (D) / Calidad (I) / Volumen (J) / Precio minimo (K) /
# Precio maximo (L) / Precio promedio ponderado (M) / Origen (O) /
# Precio $/Kg (P) values refreshed to this week's figures.

$ws.Range("D2").Value = 44159
$ws.Range("J2").Value = 1100
$ws.Range("D3").Value = 44159
$ws.Range("J3").Value = 800
$ws.Range("D4").Value = 44169
$ws.Range("J4").Value = 950
$ws.Range("D5").Value = 44169
$ws.Range("D6").Value = 44179
$ws.Range("J6").Value = 980
$ws.Range("K6").Value = 1200
$ws.Range("L6").Value = 1200
$ws.Range("M6").Value = 1200
$ws.Range("O6").Value = "Región Metropolitana"
$ws.Range("P6").Value = 1200
$ws.Range("D7").Value = 44172
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 600
$ws.Range("K7").Value = 1300
$ws.Range("L7").Value = 1300
$ws.Range("M7").Value = 1300
$ws.Range("P7").Value = 1300
$ws.Range("D8").Value = 44172
$ws.Range("I8").Value = "Segunda"
$ws.Range("J8").Value = 550
$ws.Range("K8").Value = 1000
$ws.Range("L8").Value = 1000
$ws.Range("M8").Value = 1000
$ws.Range("P8").Value = 1000
$ws.Range("D9").Value = 44162
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 1200
$ws.Range("K9").Value = 1300
$ws.Range("L9").Value = 1300
$ws.Range("M9").Value = 1300
$ws.Range("P9").Value = 1300
$ws.Range("D10").Value = 44162
$ws.Range("I10").Value = "Segunda"
$ws.Range("J10").Value = 800
$ws.Range("K10").Value = 1000
$ws.Range("L10").Value = 1000
$ws.Range("M10").Value = 1000
$ws.Range("P10").Value = 1000
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 2800
$ws.Range("K11").Value = 1200
$ws.Range("L11").Value = 1250
$ws.Range("M11").Value = 1221
$ws.Range("P11").Value = 1221
$ws.Range("D12").Value = 44174
$ws.Range("I12").Value = "Segunda"
$ws.Range("J12").Value = 1300
$ws.Range("K12").Value = 1000
$ws.Range("L12").Value = 1000
$ws.Range("M12").Value = 1000
$ws.Range("P12").Value = 1000
$ws.Range("D13").Value = 44181
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 1000
$ws.Range("K13").Value = 1300
$ws.Range("L13").Value = 1300
$ws.Range("M13").Value = 1300
$ws.Range("P13").Value = 1300
$ws.Range("D14").Value = 44181
$ws.Range("I14").Value = "Segunda"
$ws.Range("J14").Value = 900
$ws.Range("K14").Value = 900
$ws.Range("L14").Value = 900
$ws.Range("M14").Value = 900
$ws.Range("O14").Value = "Provincia de Quillota"
$ws.Range("P14").Value = 900
$ws.Range("D15").Value = 44165
$ws.Range("J15").Value = 720
$ws.Range("L15").Value = 1200
$ws.Range("M15").Value = 1200
$ws.Range("P15").Value = 1200
$ws.Range("D16").Value = 44165
$ws.Range("J16").Value = 750
$ws.Range("D17").Value = 44176
$ws.Range("J17").Value = 2500
$ws.Range("L17").Value = 1300
$ws.Range("M17").Value = 1256
$ws.Range("P17").Value = 1256
$ws.Range("D18").Value = 44176
$ws.Range("J18").Value = 1500
$ws.Range("D19").Value = 44168
$ws.Range("J19").Value = 1200
$ws.Range("D20").Value = 44168
$ws.Range("J20").Value = 850
$ws.Range("D21").Value = 44167
$ws.Range("J21").Value = 1430
$ws.Range("K21").Value = 1200
$ws.Range("M21").Value = 1248
$ws.Range("P21").Value = 1248
$ws.Range("D22").Value = 44167
$ws.Range("J22").Value = 350
$ws.Range("D23").Value = 44161
$ws.Range("J23").Value = 1600
$ws.Range("D24").Value = 44161
$ws.Range("J24").Value = 1850
$ws.Range("D25").Value = 44160
$ws.Range("J25").Value = 750
$ws.Range("K25").Value = 1300
$ws.Range("M25").Value = 1300
$ws.Range("P25").Value = 1300
$ws.Range("D26").Value = 44160
$ws.Range("J26").Value = 850
$ws.Range("D27").Value = 44175
$ws.Range("J27").Value = 1500
$ws.Range("D28").Value = 44175
$ws.Range("J28").Value = 1450
$ws.Range("K28").Value = 1000
$ws.Range("L28").Value = 1000
$ws.Range("M28").Value = 1000
$ws.Range("P28").Value = 1000
